$d = $word.ActiveDocument

# "June 18, 2022" -> "June 19, 2022" (Replace=2/wdReplaceAll covers every
# occurrence of the substring across the document, including the ones
# embedded in " on June 18, 2022." and " license is suspended from June 18, 2022")
$d.Content.Find.Execute("June 18, 2022", $true, $false, $false, $false, $false,
                         $true, 1, $false, "June 19, 2022", 2)

# "August 17, 2022" -> "August 18, 2022"
$d.Content.Find.Execute("August 17, 2022", $true, $false, $false, $false, $false,
                         $true, 1, $false, "August 18, 2022", 2)
